$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44231
$ws.Range("J2").Value = 1300
$ws.Range("K2").Value = 450
$ws.Range("L2").Value = 480
$ws.Range("M2").Value = 465
$ws.Range("P2").Value = 465
$ws.Range("D3").Value = 44428
$ws.Range("I3").Value = '1a nueva(o)'
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 580
$ws.Range("L3").Value = 600
$ws.Range("M3").Value = 590
$ws.Range("O3").Value = 'Perú'
$ws.Range("P3").Value = 590
$ws.Range("D4").Value = 44490
$ws.Range("H4").Value = 'Camote'
$ws.Range("K4").Value = 450
$ws.Range("L4").Value = 480
$ws.Range("M4").Value = 465
$ws.Range("P4").Value = 465
$ws.Range("D5").Value = 44476
$ws.Range("J5").Value = 1200
$ws.Range("K5").Value = 480
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = 490
$ws.Range("P5").Value = 490
$ws.Range("D6").Value = 44243
$ws.Range("I6").Value = '2a nueva(o)'
$ws.Range("J6").Value = 1600
$ws.Range("K6").Value = 450
$ws.Range("L6").Value = 480
$ws.Range("M6").Value = 465
$ws.Range("P6").Value = 465
$ws.Range("D7").Value = 44179
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 1350
$ws.Range("L7").Value = 1400
$ws.Range("M7").Value = 1375
$ws.Range("P7").Value = 1375
$ws.Range("D8").Value = 44399
$ws.Range("I8").Value = '1a (guarda)'
$ws.Range("J8").Value = 800
$ws.Range("K8").Value = 450
$ws.Range("L8").Value = 480
$ws.Range("M8").Value = 465
$ws.Range("O8").Value = 'Provincia de Melipilla'
$ws.Range("P8").Value = 465
$ws.Range("D9").Value = 44175
$ws.Range("K9").Value = 1500
$ws.Range("L9").Value = 1550
$ws.Range("M9").Value = 1525
$ws.Range("P9").Value = 1525
$ws.Range("D10").Value = 44168
$ws.Range("I10").Value = '1a nueva(o)'
$ws.Range("J10").Value = 1200
$ws.Range("K10").Value = 1500
$ws.Range("L10").Value = 1700
$ws.Range("M10").Value = 1600
$ws.Range("O10").Value = 'Perú'
$ws.Range("P10").Value = 1600
$ws.Range("D11").Value = 44483
$ws.Range("J11").Value = 1300
$ws.Range("K11").Value = 550
$ws.Range("L11").Value = 580
$ws.Range("M11").Value = 565
$ws.Range("P11").Value = 565
$ws.Range("D12").Value = 44349
$ws.Range("H12").Value = 'Pachia'
$ws.Range("K12").Value = 730
$ws.Range("L12").Value = 750
$ws.Range("M12").Value = 740
$ws.Range("P12").Value = 740
$ws.Range("D13").Value = 44211
$ws.Range("J13").Value = 1600
$ws.Range("K13").Value = 500
$ws.Range("L13").Value = 550
$ws.Range("M13").Value = 525
$ws.Range("O13").Value = 'Región de O''Higgins'
$ws.Range("P13").Value = 525
$ws.Range("D15").Value = 44469
$ws.Range("J15").Value = 1200
$ws.Range("K15").Value = 600
$ws.Range("L15").Value = 650
$ws.Range("M15").Value = 625
$ws.Range("O15").Value = 'Perú'
$ws.Range("P15").Value = 625
$ws.Range("D16").Value = 44453
$ws.Range("J16").Value = 800
$ws.Range("K16").Value = 630
$ws.Range("M16").Value = 640
$ws.Range("P16").Value = 640
$ws.Range("D17").Value = 44301
$ws.Range("I17").Value = '2a nueva(o)'
$ws.Range("K17").Value = 400
$ws.Range("L17").Value = 430
$ws.Range("M17").Value = 415
$ws.Range("O17").Value = 'Provincia de Melipilla'
$ws.Range("P17").Value = 415
$ws.Range("D18").Value = 44322
$ws.Range("I18").Value = '1a (cosecha)'
$ws.Range("K18").Value = 350
$ws.Range("L18").Value = 400
$ws.Range("M18").Value = 375
$ws.Range("O18").Value = 'Región del Maule'
$ws.Range("P18").Value = 375
$ws.Range("D19").Value = 44238
$ws.Range("I19").Value = '1a nueva(o)'
$ws.Range("J19").Value = 1250
$ws.Range("K19").Value = 430
$ws.Range("L19").Value = 450
$ws.Range("M19").Value = 440
$ws.Range("O19").Value = 'Perú'
$ws.Range("P19").Value = 440
$ws.Range("D20").Value = 44201
$ws.Range("J20").Value = 1360
$ws.Range("K20").Value = 730
$ws.Range("L20").Value = 750
$ws.Range("M20").Value = 740
$ws.Range("P20").Value = 740
